$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Ruoran%Wang%NULL%1,                 Min%He%NULL%1,                 Wanhong%Yin%NULL%1,                 Xuelian%Liao%NULL%1,                 Bo%Wang%NULL%6,                 Xiaodong%Jin%NULL%1,                 Yao%Ma%NULL%2,                 Jirong%Yue%NULL%1,                 Lang%Bai%NULL%1,                 Dan%Liu%NULL%3,                 Ting%Zhu%NULL%1,                 Zhixin%Huang%NULL%1,                 Yan%Kang%kangyan@scu.edu.cn%2]'
$ws.Range("E3").Value = '[Gang%Wu%NULL%1,                 Shuchang%Zhou%NULL%1,                 Yujin%Wang%NULL%1,                 Wenzhi%Lv%NULL%2,                 Shili%Wang%NULL%1,                 Ting%Wang%751884926@qq.com%1,                 Xiaoming%Li%lilyboston2002@qq.com%0]'
$ws.Range("E4").Value = '[Arjun S%Yadaw%NULL%1,                 Yan-chak%Li%NULL%1,                 Sonali%Bose%NULL%1,                 Ravi%Iyengar%NULL%1,                 Supinda%Bunyavanich%NULL%1,                 Gaurav%Pandey%NULL%1]'
$ws.Range("E5").Value = '[ Li%Yan%null%1,        Hai-Tao%Zhang%null%1,        Jorge%Goncalves%null%1,        Yang%Xiao%null%1,        Maolin%Wang%null%1,        Yuqi%Guo%null%1,        Chuan%Sun%null%1,        Xiuchuan%Tang%null%1,        Liang%Jing%null%1,        Mingyang%Zhang%null%1,        Xiang%Huang%null%1,        Ying%Xiao%null%1,        Haosen%Cao%null%1,        Yanyan%Chen%null%1,        Tongxin%Ren%null%1,        Fang%Wang%null%1,        Yaru%Xiao%null%1,        Sufang%Huang%null%1,        Xi%Tan%null%1,        Niannian%Huang%null%1,        Bo%Jiao%null%1,        Cheng%Cheng%null%1,        Yong%Zhang%null%1,        Ailin%Luo%null%1,        Laurent%Mombaerts%null%1,        Junyang%Jin%null%1,        Zhiguo%Cao%null%1,        Shusheng%Li%null%1,        Hui%Xu%null%1,        Ye%Yuan%null%1]'
$ws.Range("E6").Value = '[Shuai%Zhang%NULL%1,                 Mengfei%Guo%NULL%1,                 Limin%Duan%NULL%1,                 Feng%Wu%NULL%1,                 Guorong%Hu%NULL%1,                 Zhihui%Wang%NULL%1,                 Qi%Huang%NULL%1,                 Tingting%Liao%NULL%1,                 Juanjuan%Xu%NULL%1,                 Yanling%Ma%NULL%1,                 Zhilei%Lv%NULL%1,                 Wenjing%Xiao%NULL%1,                 Zilin%Zhao%NULL%1,                 Xueyun%Tan%NULL%1,                 Daquan%Meng%NULL%1,                 Shujing%Zhang%NULL%1,                 E%Zhou%NULL%1,                 Zhengrong%Yin%NULL%1,                 Wei%Geng%NULL%1,                 Xuan%Wang%NULL%1,                 Jianchu%Zhang%NULL%1,                 Jianguo%Chen%Chenj@mails.tjmu.edu.cn%1,                 Yu%Zhang%whxhzy@163.com%1,                 Yang%Jin%whuhjy@126.com%1]'
$ws.Range("E7").Value = '[Angelo%Zinellu%azinellu@uniss.it%1,                 Francesco%Arru%NULL%2,                 Francesco%Arru%NULL%0,                 Andrea%De Vito%NULL%3,                 Alessandro%Sassu%NULL%2,                 Alessandro%Sassu%NULL%0,                 Giovanni%Valdes%NULL%1,                 Valentina%Scano%NULL%1,                 Elisabetta%Zinellu%NULL%1,                 Roberto%Perra%NULL%1,                 Giordano%Madeddu%NULL%3,                 Ciriaco%Carru%NULL%1,                 Pietro%Pirina%NULL%3,                 Arduino A.%Mangoni%NULL%2,                 Arduino A.%Mangoni%NULL%0,                 Sergio%Babudieri%NULL%4,                 Sergio%Babudieri%NULL%0,                 Alessandro G.%Fois%NULL%1]'
$ws.Range("E8").Value = '[Juan%Torres-Macho%NULL%1,                 Pablo%Ryan%NULL%1,                 Jorge%Valencia%NULL%2,                 Jorge%Valencia%NULL%0,                 Mario%Pérez-Butragueño%NULL%2,                 Mario%Pérez-Butragueño%NULL%0,                 Eva%Jiménez%NULL%2,                 Eva%Jiménez%NULL%0,                 Mario%Fontán-Vela%NULL%1,                 Elsa%Izquierdo-García%NULL%2,                 Elsa%Izquierdo-García%NULL%0,                 Inés%Fernandez-Jimenez%NULL%1,                 Elena%Álvaro-Alonso%NULL%1,                 Andrea%Lazaro%NULL%2,                 Andrea%Lazaro%NULL%0,                 Marta%Alvarado%NULL%1,                 Helena%Notario%NULL%1,                 Salvador%Resino%NULL%1,                 Daniel%Velez-Serrano%NULL%1,                 Alejandro%Meca%NULL%1]'
$ws.Range("E9").Value = '[Qingquan%Liu%NULL%1,                 Yiru%Wang%NULL%2,                 Xuecheng%Zhao%NULL%1,                 Lixuan%Wang%NULL%1,                 Feng%Liu%NULL%3,                 Tao%Wang%NULL%15,                 Dawei%Ye%NULL%2,                 Yongman%Lv%NULL%2]'
$ws.Range("E10").Value = '[Bianca%Magro%NULL%1,                 Valentina%Zuccaro%NULL%1,                 Luca%Novelli%NULL%1,                 Lorenzo%Zileri%NULL%2,                 Lorenzo%Zileri%NULL%0,                 Ciro%Celsa%NULL%1,                 Federico%Raimondi%NULL%2,                 Federico%Raimondi%NULL%0,                 Mauro%Gori%NULL%1,                 Giulia%Cammà%NULL%1,                 Salvatore%Battaglia%NULL%1,                 Vincenzo Giuseppe%Genova%NULL%1,                 Laura%Paris%NULL%1,                 Matteo%Tacelli%NULL%1,                 Francesco Antonio%Mancarella%NULL%1,                 Marco%Enea%NULL%1,                 Massimo%Attanasio%NULL%1,                 Michele%Senni%NULL%1,                 Fabiano%Di Marco%NULL%1,                 Luca Ferdinando%Lorini%NULL%1,                 Stefano%Fagiuoli%NULL%1,                 Raffaele%Bruno%NULL%2,                 Calogero%Cammà%NULL%1,                 Antonio%Gasbarrini%NULL%3,                 Francesco%Di Gennaro%NULL%2,                 Francesco%Di Gennaro%NULL%0]'
$ws.Range("E11").Value = '[Ze-yang%Ding%NULL%1,                 Gan-xun%Li%NULL%1,                 Lin%Chen%NULL%3,                 Chang%Shu%NULL%2,                 Jia%Song%NULL%1,                 Wei%Wang%NULL%2,                 Yu-wei%Wang%NULL%1,                 Qian%Chen%NULL%2,                 Guan-nan%Jin%NULL%1,                 Tong-tong%Liu%NULL%1,                 Jun-nan%Liang%NULL%1,                 Peng%Zhu%NULL%1,                 Wei%Zhu%NULL%5,                 Yong%Li%NULL%2,                 Bin-hao%Zhang%NULL%1,                 Huan%Feng%NULL%1,                 Wan-guang%Zhang%NULL%1,                 Zhen-yu%Yin%NULL%1,                 Wen-kui%Yu%NULL%1,                 Yang%Yang%NULL%2,                 Hua-qiu%Zhang%NULL%1,                 Zhou-ping%Tang%NULL%1,                 Hui%Wang%NULL%4,                 Jun-bo%Hu%NULL%1,                 Ji-hong%Liu%NULL%1,                 Ping%Yin%NULL%1,                 Xiao-ping%Chen%NULL%1,                 Bixiang%Zhang%NULL%1,                 NULL%NULL%NULL%0]'
$ws.Range("E12").Value = '[Hazal Cansu%Acar%hazal.acar@istanbul.edu.tr%1,                 Günay%Can%NULL%2,                 Günay%Can%NULL%0,                 Rıdvan%Karaali%NULL%1,                 Şermin%Börekçi%NULL%1,                 İlker İnanç%Balkan%NULL%1,                 Bilun%Gemicioğlu%NULL%1,                 Dildar%Konukoğlu%NULL%1,                 Ethem%Erginöz%NULL%1,                 Mehmet Sarper%Erdoğan%NULL%1,                 Fehmi%Tabak%NULL%1]'
$ws.Range("E13").Value = '[Mattia%Bellan%NULL%1,                 Danila%Azzolina%NULL%1,                 Eyal%Hayden%NULL%2,                 Eyal%Hayden%NULL%0,                 Gianluca%Gaidano%NULL%2,                 Gianluca%Gaidano%NULL%0,                 Mario%Pirisi%NULL%2,                 Mario%Pirisi%NULL%0,                 Antonio%Acquaviva%NULL%2,                 Antonio%Acquaviva%NULL%0,                 Gianluca%Aimaretti%NULL%2,                 Gianluca%Aimaretti%NULL%0,                 Paolo%Aluffi Valletti%NULL%3,                 Paolo%Aluffi Valletti%NULL%0,                 Roberto%Angilletta%NULL%2,                 Roberto%Angilletta%NULL%0,                 Roberto%Arioli%NULL%2,                 Roberto%Arioli%NULL%0,                 Gian Carlo%Avanzi%NULL%2,                 Gian Carlo%Avanzi%NULL%0,                 Gianluca%Avino%NULL%2,                 Gianluca%Avino%NULL%0,                 Piero Emilio%Balbo%NULL%2,                 Piero Emilio%Balbo%NULL%0,                 Giulia%Baldon%NULL%2,                 Giulia%Baldon%NULL%0,                 Francesca%Baorda%NULL%2,                 Francesca%Baorda%NULL%0,                 Emanuela%Barbero%NULL%2,                 Emanuela%Barbero%NULL%0,                 Alessio%Baricich%NULL%2,                 Alessio%Baricich%NULL%0,                 Michela%Barini%NULL%2,                 Michela%Barini%NULL%0,                 Francesco%Barone-Adesi%NULL%2,                 Francesco%Barone-Adesi%NULL%0,                 Sofia%Battistini%NULL%2,                 Sofia%Battistini%NULL%0,                 Michela%Beltrame%NULL%2,                 Michela%Beltrame%NULL%0,                 Matteo%Bertoli%NULL%2,                 Matteo%Bertoli%NULL%0,                 Stephanie%Bertolin%NULL%2,                 Stephanie%Bertolin%NULL%0,                 Marinella%Bertolotti%NULL%2,                 Marinella%Bertolotti%NULL%0,                 Marta%Betti%NULL%2,                 Marta%Betti%NULL%0,                 Flavio%Bobbio%NULL%2,                 Flavio%Bobbio%NULL%0,                 Paolo%Boffano%NULL%2,                 Paolo%Boffano%NULL%0,                 Lucio%Boglione%NULL%2,                 Lucio%Boglione%NULL%0,                 Silvio%Borrè%NULL%2,                 Silvio%Borrè%NULL%0,                 Matteo%Brucoli%NULL%2,                 Matteo%Brucoli%NULL%0,                 Elisa%Calzaducca%NULL%2,                 Elisa%Calzaducca%NULL%0,                 Edoardo%Cammarata%NULL%2,                 Edoardo%Cammarata%NULL%0,                 Vincenzo%Cantaluppi%NULL%2,                 Vincenzo%Cantaluppi%NULL%0,                 Roberto%Cantello%NULL%2,                 Roberto%Cantello%NULL%0,                 Andrea%Capponi%NULL%2,                 Andrea%Capponi%NULL%0,                 Alessandro%Carriero%NULL%2,                 Alessandro%Carriero%NULL%0,                 Giuseppe Francesco%Casciaro%NULL%2,                 Giuseppe Francesco%Casciaro%NULL%0,                 Luigi Mario%Castello%NULL%2,                 Luigi Mario%Castello%NULL%0,                 Federico%Ceruti%NULL%2,                 Federico%Ceruti%NULL%0,                 Guido%Chichino%NULL%2,                 Guido%Chichino%NULL%0,                 Emilio%Chirico%NULL%2,                 Emilio%Chirico%NULL%0,                 Carlo%Cisari%NULL%1,                 Micol Giulia%Cittone%NULL%2,                 Micol Giulia%Cittone%NULL%0,                 Crizia%Colombo%NULL%2,                 Crizia%Colombo%NULL%0,                 Cristoforo%Comi%NULL%2,                 Cristoforo%Comi%NULL%0,                 Eleonora%Croce%NULL%2,                 Eleonora%Croce%NULL%0,                 Tommaso%Daffara%NULL%2,                 Tommaso%Daffara%NULL%0,                 Pietro%Danna%NULL%2,                 Pietro%Danna%NULL%0,                 Francesco%Della Corte%NULL%2,                 Francesco%Della Corte%NULL%0,                 Simona%De Vecchi%NULL%2,                 Simona%De Vecchi%NULL%0,                 Umberto%Dianzani%NULL%2,                 Umberto%Dianzani%NULL%0,                 Davide%Di Benedetto%NULL%2,                 Davide%Di Benedetto%NULL%0,                 Elia%Esposto%NULL%2,                 Elia%Esposto%NULL%0,                 Fabrizio%Faggiano%NULL%2,                 Fabrizio%Faggiano%NULL%0,                 Zeno%Falaschi%NULL%2,                 Zeno%Falaschi%NULL%0,                 Daniela%Ferrante%NULL%2,                 Daniela%Ferrante%NULL%0,                 Alice%Ferrero%NULL%2,                 Alice%Ferrero%NULL%0,                 Ileana%Gagliardi%NULL%2,                 Ileana%Gagliardi%NULL%0,                 Alessandra%Galbiati%NULL%2,                 Alessandra%Galbiati%NULL%0,                 Silvia%Gallo%NULL%2,                 Silvia%Gallo%NULL%0,                 Pietro Luigi%Garavelli%NULL%2,                 Pietro Luigi%Garavelli%NULL%0,                 Clara Ada%Gardino%NULL%2,                 Clara Ada%Gardino%NULL%0,                 Massimiliano%Garzaro%NULL%3,                 Massimiliano%Garzaro%NULL%0,                 Maria Luisa%Gastaldello%NULL%2,                 Maria Luisa%Gastaldello%NULL%0,                 Francesco%Gavelli%NULL%2,                 Francesco%Gavelli%NULL%0,                 Alessandra%Gennari%NULL%2,                 Alessandra%Gennari%NULL%0,                 Greta Maria%Giacomini%NULL%2,                 Greta Maria%Giacomini%NULL%0,                 Irene%Giacone%NULL%2,                 Irene%Giacone%NULL%0,                 Valentina%Giai Via%NULL%2,                 Valentina%Giai Via%NULL%0,                 Francesca%Giolitti%NULL%2,                 Francesca%Giolitti%NULL%0,                 Laura Cristina%Gironi%NULL%2,                 Laura Cristina%Gironi%NULL%0,                 Carla%Gramaglia%NULL%2,                 Carla%Gramaglia%NULL%0,                 Leonardo%Grisafi%NULL%2,                 Leonardo%Grisafi%NULL%0,                 Ilaria%Inserra%NULL%2,                 Ilaria%Inserra%NULL%0,                 Marco%Invernizzi%NULL%2,                 Marco%Invernizzi%NULL%0,                 Marco%Krengli%NULL%2,                 Marco%Krengli%NULL%0,                 Emanuela%Labella%NULL%2,                 Emanuela%Labella%NULL%0,                 Irene Cecilia%Landi%NULL%2,                 Irene Cecilia%Landi%NULL%0,                 Raffaella%Landi%NULL%2,                 Raffaella%Landi%NULL%0,                 Ilaria%Leone%NULL%2,                 Ilaria%Leone%NULL%0,                 Veronica%Lio%NULL%2,                 Veronica%Lio%NULL%0,                 Luca%Lorenzini%NULL%2,                 Luca%Lorenzini%NULL%0,                 Antonio%Maconi%NULL%2,                 Antonio%Maconi%NULL%0,                 Mario%Malerba%NULL%2,                 Mario%Malerba%NULL%0,                 Giulia Francesca%Manfredi%NULL%2,                 Giulia Francesca%Manfredi%NULL%0,                 Maria%Martelli%NULL%2,                 Maria%Martelli%NULL%0,                 Letizia%Marzari%NULL%2,                 Letizia%Marzari%NULL%0,                 Paolo%Marzullo%NULL%2,                 Paolo%Marzullo%NULL%0,                 Marco%Mennuni%NULL%2,                 Marco%Mennuni%NULL%0,                 Claudia%Montabone%NULL%2,                 Claudia%Montabone%NULL%0,                 Umberto%Morosini%NULL%2,                 Umberto%Morosini%NULL%0,                 Marco%Mussa%NULL%2,                 Marco%Mussa%NULL%0,                 Ilaria%Nerici%NULL%2,                 Ilaria%Nerici%NULL%0,                 Alessandro%Nuzzo%NULL%2,                 Alessandro%Nuzzo%NULL%0,                 Carlo%Olivieri%NULL%2,                 Carlo%Olivieri%NULL%0,                 Samuel Alberto%Padelli%NULL%2,                 Samuel Alberto%Padelli%NULL%0,                 Massimiliano%Panella%NULL%2,                 Massimiliano%Panella%NULL%0,                 Andrea%Parisini%NULL%2,                 Andrea%Parisini%NULL%0,                 Alessio%Paschè%NULL%2,                 Alessio%Paschè%NULL%0,                 Filippo%Patrucco%NULL%2,                 Filippo%Patrucco%NULL%0,                 Giuseppe%Patti%NULL%2,                 Giuseppe%Patti%NULL%0,                 Alberto%Pau%NULL%2,                 Alberto%Pau%NULL%0,                 Anita Rebecca%Pedrinelli%NULL%2,                 Anita Rebecca%Pedrinelli%NULL%0,                 Ilaria%Percivale%NULL%2,                 Ilaria%Percivale%NULL%0,                 Luca%Ragazzoni%NULL%2,                 Luca%Ragazzoni%NULL%0,                 Roberta%Re%NULL%2,                 Roberta%Re%NULL%0,                 Cristina%Rigamonti%NULL%2,                 Cristina%Rigamonti%NULL%0,                 Eleonora%Rizzi%NULL%2,                 Eleonora%Rizzi%NULL%0,                 Andrea%Rognoni%NULL%2,                 Andrea%Rognoni%NULL%0,                 Annalisa%Roveta%NULL%2,                 Annalisa%Roveta%NULL%0,                 Luigia%Salamina%NULL%2,                 Luigia%Salamina%NULL%0,                 Matteo%Santagostino%NULL%2,                 Matteo%Santagostino%NULL%0,                 Massimo%Saraceno%NULL%2,                 Massimo%Saraceno%NULL%0,                 Paola%Savoia%NULL%2,                 Paola%Savoia%NULL%0,                 Marco%Sciarra%NULL%2,                 Marco%Sciarra%NULL%0,                 Andrea%Schimmenti%NULL%2,                 Andrea%Schimmenti%NULL%0,                 Lorenza%Scotti%NULL%2,                 Lorenza%Scotti%NULL%0,                 Enrico%Spinoni%NULL%2,                 Enrico%Spinoni%NULL%0,                 Carlo%Smirne%NULL%2,                 Carlo%Smirne%NULL%0,                 Vanessa%Tarantino%NULL%2,                 Vanessa%Tarantino%NULL%0,                 Paolo Amedeo%Tillio%NULL%2,                 Paolo Amedeo%Tillio%NULL%0,                 Stelvio%Tonello%NULL%2,                 Stelvio%Tonello%NULL%0,                 Rosanna%Vaschetto%NULL%2,                 Rosanna%Vaschetto%NULL%0,                 Veronica%Vassia%NULL%2,                 Veronica%Vassia%NULL%0,                 Domenico%Zagaria%NULL%2,                 Domenico%Zagaria%NULL%0,                 Elisa%Zavattaro%NULL%2,                 Elisa%Zavattaro%NULL%0,                 Patrizia%Zeppegno%NULL%2,                 Patrizia%Zeppegno%NULL%0,                 Francesca%Zottarelli%NULL%2,                 Francesca%Zottarelli%NULL%0,                 Pier Paolo%Sainaghi%pierpaolo.sainaghi@med.uniupo.it%2,                 Pier Paolo%Sainaghi%pierpaolo.sainaghi@med.uniupo.it%0]'
$ws.Range("E14").Value = '[Giulia%Besutti%giulia.besutti@ausl.re.it%1,                 Marta%Ottone%NULL%2,                 Marta%Ottone%NULL%0,                 Tommaso%Fasano%NULL%1,                 Pierpaolo%Pattacini%NULL%1,                 Valentina%Iotti%NULL%1,                 Lucia%Spaggiari%NULL%1,                 Riccardo%Bonacini%NULL%1,                 Andrea%Nitrosi%NULL%1,                 Efrem%Bonelli%NULL%1,                 Simone%Canovi%NULL%1,                 Rossana%Colla%NULL%1,                 Alessandro%Zerbini%NULL%1,                 Marco%Massari%NULL%1,                 Ivana%Lattuada%NULL%1,                 Anna Maria%Ferrari%NULL%1,                 Paolo%Giorgi Rossi%NULL%1,                 Massimo%Costantini%NULL%1,                 Roberto%Grilli%NULL%1,                 Massimiliano%Marino%NULL%1,                 Giulio%Formoso%NULL%1,                 Debora%Formisano%NULL%1,                 Emanuela%Bedeschi%NULL%1,                 Cinzia%Perilli%NULL%1,                 Elisabetta%La Rosa%NULL%1,                 Eufemia%Bisaccia%NULL%1,                 Ivano%Venturi%NULL%1,                 Massimo%Vicentini%NULL%1,                 Cinzia%Campari%NULL%1,                 Francesco%Gioia%NULL%1,                 Serena%Broccoli%NULL%1,                 Pamela%Mancuso%NULL%1,                 Marco%Foracchia%NULL%1,                 Mirco%Pinotti%NULL%1,                 Nicola%Facciolongo%NULL%1,                 Laura%Trabucco%NULL%1,                 Stefano%De Pietri%NULL%1,                 Giorgio Francesco%Danelli%NULL%1,                 Laura%Albertazzi%NULL%1,                 Enrica%Bellesia%NULL%1,                 Mattia%Corradini%NULL%1,                 Elena%Magnani%NULL%1,                 Annalisa%Pilia%NULL%1,                 Alessandra%Polese%NULL%1,                 Silvia Storchi%Incerti%NULL%1,                 Piera%Zaldini%NULL%1,                 Bonanno%Orsola%NULL%1,                 Matteo%Revelli%NULL%1,                 Carlo%Salvarani%NULL%1,                 Carmine%Pinto%NULL%1,                 Francesco%Venturelli%NULL%1]'
$ws.Range("E15").Value = '[Bo%Chen%NULL%1,                 Hong-Qiu%Gu%NULL%1,                 Yi%Liu (刘艺)%NULL%1,                 Guqin%Zhang%NULL%1,                 Hang%Yang%NULL%1,                 Huifang%Hu%NULL%1,                 Chenyang%Lu%NULL%1,                 Yang%Li%NULL%3,                 Liyi%Wang%NULL%1,                 Yi%Liu (刘毅)%yi2006liu@163.com%1,                 Yi%Zhao%zhao.y1977@163.com%1,                 Huaqin%Pan%phq2012@whu.edu.cn%1]'
$ws.Range("E16").Value = '[Muhammad E. H.%Chowdhury%mchowdhury@qu.edu.qa%1,                 Tawsifur%Rahman%NULL%2,                 Tawsifur%Rahman%NULL%0,                 Amith%Khandakar%NULL%3,                 Somaya%Al-Madeed%NULL%2,                 Susu M.%Zughaier%NULL%5,                 Suhail A. R.%Doi%NULL%3,                 Hanadi%Hassen%NULL%1,                 Mohammad T.%Islam%NULL%1]'
$ws.Range("E17").Value = '[Yi-Min%Dong%NULL%1,                 Jia%Sun%NULL%1,                 Yi-Xin%Li%NULL%1,                 Qian%Chen%NULL%0,                 Qing-Quan%Liu%NULL%1,                 Zhou%Sun%NULL%1,                 Ran%Pang%NULL%1,                 Fei%Chen%NULL%1,                 Bing-Yang%Xu%NULL%1,                 Anne%Manyande%NULL%1,                 Taane G%Clark%NULL%1,                 Jin-Ping%Li%NULL%1,                 Ilkay Erdogan%Orhan%NULL%1,                 Yu-Ke%Tian%NULL%1,                 Tao%Wang%wt7636@126.com%0,                 Wei%Wu%wt7636@126.com%1,                 Da-Wei%Ye%wt7636@126.com%1]'
$ws.Range("E18").Value = '[Mark H.%Ebell%xref no email%1,          Xinyan%Cai%xref no email%1,          Robert%Lennon%xref no email%1,          Derjung M.%Tarn%xref no email%1,          Arch G.%Mainous%xref no email%1,          Aleksandra E.%Zgierska%xref no email%1,          Bruce%Barrett%xref no email%1,          Wen-Jan%Tuan%xref no email%1,          Kevin%Maloy%xref no email%1,          Munish%Goyal%xref no email%1,          Alex%Krist%xref no email%1]'
$ws.Range("E19").Value = '[Cao%Y%coreGivesNoEmail%4,              Imam%Z%coreGivesNoEmail%1,              Lippi%G%coreGivesNoEmail%1,              Oran%DP%coreGivesNoEmail%1,              Shi%S%coreGivesNoEmail%1]'
$ws.Range("E20").Value = '[Jialin%He%NULL%1,                 Caiping%Song%NULL%1,                 En%Liu%NULL%1,                 Xi%Liu%NULL%1,                 Hao%Wu%NULL%1,                 Hui%Lin%NULL%1,                 Yuliang%Liu%NULL%1,                 Qi%Li%NULL%1,                 Zhi%Xu%NULL%1,                 XiaoBao%Ren%NULL%1,                 Cheng%Zhang%NULL%1,                 Wenjing%Zhang%NULL%1,                 Wei%Duan%NULL%2,                 Yongfeng%Tian%NULL%1,                 Ping%Li%NULL%1,                 Mingdong%Hu%NULL%1,                 Shiming%Yang%NULL%1,                 Yu%Xu%NULL%1]'
$ws.Range("E21").Value = '[Meng%Jiang%NULL%1,                 Changli%Li%NULL%1,                 Li%Zheng%NULL%1,                 Wenzhi%Lv%NULL%0,                 Zhigang%He%NULL%1,                 Xinwu%Cui%NULL%1,                 Christoph F.%Dietrich%NULL%1]'
$ws.Range("E22").Value = '[Zan%Ke%xref no email%1,          Liang%Li%xref no email%1,          Li%Wang%xref no email%1,          Huan%Liu%xref no email%2,          Xuefang%Lu%xref no email%1,          Feifei%Zeng%xref no email%1,          Yunfei%Zha%xref no email%1]'
$ws.Range("E23").Value = '[Do Hyoung%Kim%NULL%1,                 Hayne Cho%Park%NULL%2,                 Hayne Cho%Park%NULL%0,                 Ajin%Cho%NULL%1,                 Juhee%Kim%NULL%1,                 Kyu-sang%Yun%NULL%1,                 Jinseog%Kim%NULL%1,                 Young-Ki%Lee%NULL%1,                 Sinan%Kardes.%NULL%2,                 Sinan%Kardes.%NULL%0]'
$ws.Range("E24").Value = '[Jian%Li%xref no email%1,          Luyu%Yang%xref no email%1,          Qian%Zeng%xref no email%1,          Qingyun%Li%xref no email%1,          Zhitao%Yang%xref no email%1,          Lizhong%Han%xref no email%1,          Xiaodong%Huang%xref no email%1,          Erzhen%Chen%xref no email%1]'
$ws.Range("E25").Value = '[Li%Li%NULL%1,                 Xiaoyu%Fang%NULL%1,                 Lixia%Cheng%NULL%1,                 Penghao%Wang%NULL%1,                 Shen%Li%NULL%1,                 Hao%Yu%NULL%1,                 Yao%Zhang%NULL%2,                 Nan%Jiang%NULL%1,                 Tingting%Zeng%NULL%1,                 Chao%Hou%NULL%1,                 Jing%Zhou%NULL%0,                 Shiru%Li%NULL%1,                 Yingzi%Pan%NULL%1,                 Yitong%Li%NULL%1,                 Lili%Nie%NULL%1,                 Yang%Li%NULL%0,                 Qidi%Sun%NULL%1,                 Hong%Jia%NULL%1,                 Mengxia%Li%NULL%1,                 Guoqiang%Cao%NULL%1,                 Xiangyu%Ma%NULL%1]'
$ws.Range("E26").Value = '[Maarten C%Ottenhoff%NULL%1,                 Lucas A%Ramos%NULL%2,                 Lucas A%Ramos%NULL%0,                 Wouter%Potters%NULL%1,                 Marcus L F%Janssen%NULL%1,                 Deborah%Hubers%NULL%1,                 Shi%Hu%NULL%1,                 Egill A%Fridgeirsson%NULL%1,                 Dan%Piña-Fuentes%NULL%1,                 Rajat%Thomas%NULL%1,                 Iwan C C%van der Horst%NULL%1,                 Christian%Herff%NULL%2,                 Christian%Herff%NULL%0,                 Pieter%Kubben%NULL%1,                 Paul W G%Elbers%NULL%1,                 Henk A%Marquering%NULL%1,                 Max%Welling%NULL%1,                 Suat%Simsek%NULL%1,                 Martijn D%de Kruif%NULL%1,                 Tom%Dormans%NULL%1,                 Lucas M%Fleuren%NULL%1,                 Michiel%Schinkel%NULL%1,                 Peter G%Noordzij%NULL%1,                 Joop P%van den Bergh%NULL%2,                 Joop P%van den Bergh%NULL%0,                 Caroline E%Wyers%NULL%1,                 David T B%Buis%NULL%2,                 David T B%Buis%NULL%0,                 W Joost%Wiersinga%NULL%1,                 Ella H C%van den Hout%NULL%1,                 Auke C%Reidinga%NULL%1,                 Daisy%Rusch%NULL%1,                 Kim C E%Sigaloff%NULL%1,                 Renee A%Douma%NULL%1,                 Lianne%de Haan%NULL%1,                 Niels C%Gritters van den Oever%NULL%1,                 Roger J M W%Rennenberg%NULL%1,                 Guido A%van Wingen%NULL%1,                 Marcel J H%Aries%NULL%1,                 Martijn%Beudel%NULL%2,                 Martijn%Beudel%NULL%0,                 NULL%NULL%NULL%0]'
$ws.Range("E27").Value = '[Tawsifur%Rahman%NULL%2,                 Fajer A.%Al-Ishaq%NULL%2,                 Fajer A.%Al-Ishaq%NULL%0,                 Fatima S.%Al-Mohannadi%NULL%1,                 Reem S.%Mubarak%NULL%1,                 Maryam H.%Al-Hitmi%NULL%1,                 Khandaker Reajul%Islam%NULL%1,                 Amith%Khandakar%NULL%0,                 Ali Ait%Hssain%NULL%3,                 Ali Ait%Hssain%NULL%0,                 Somaya%Al-Madeed%NULL%0,                 Susu M.%Zughaier%NULL%0,                 Susu M.%Zughaier%NULL%0,                 Muhammad E. H.%Chowdhury%NULL%3,                 Muhammad E. H.%Chowdhury%NULL%0,                 Antonella%Santone%NULL%3,                 Antonella%Santone%NULL%0,                 Antonella%Santone%NULL%0,                 Emanuele%Neri%NULL%1]'
$ws.Range("E28").Value = '[Tawsifur%Rahman%NULL%0,                 Amith%Khandakar%NULL%0,                 Md Enamul%Hoque%NULL%2,                 Md Enamul%Hoque%NULL%0,                 Nabil%Ibtehaz%NULL%2,                 Nabil%Ibtehaz%NULL%0,                 Saad Bin%Kashem%NULL%2,                 Saad Bin%Kashem%NULL%0,                 Reehum%Masud%NULL%1,                 Lutfunnahar%Shampa%NULL%1,                 Mohammad Mehedi%Hasan%NULL%1,                 Mohammad Tariqul%Islam%NULL%1,                 Somaya%Al-Maadeed%NULL%2,                 Somaya%Al-Maadeed%NULL%0,                 Susu M.%Zughaier%NULL%0,                 Susu M.%Zughaier%NULL%0,                 Saif%Badran%NULL%2,                 Saif%Badran%NULL%0,                 Suhail A. R.%Doi%NULL%0,                 Suhail A. R.%Doi%NULL%0,                 Muhammad E. H.%Chowdhury%mchowdhury@qu.edu.qa%0]'
$ws.Range("E29").Value = '[Manuel%Sánchez-Montañés%NULL%1,                 Pablo%Rodríguez-Belenguer%NULL%2,                 Pablo%Rodríguez-Belenguer%NULL%0,                 Antonio J.%Serrano-López%NULL%1,                 Emilio%Soria-Olivas%NULL%2,                 Emilio%Soria-Olivas%NULL%0,                 Yasser%Alakhdar-Mohmara%NULL%1]'
$ws.Range("E30").Value = '[Hai%Wang%NULL%1,                 Haibo%Ai%NULL%1,                 Yunong%Fu%NULL%1,                 Qinglin%Li%NULL%1,                 Ruixia%Cui%NULL%1,                 Xiaohua%Ma%NULL%1,                 Yan-fen%Ma%NULL%1,                 Zi%Wang%NULL%1,                 Tong%Liu%NULL%1,                 Yunxiang%Long%NULL%1,                 Kai%Qu%NULL%1,                 Chang%Liu%NULL%1,                 Jingyao%Zhang%NULL%1]'
$ws.Range("E31").Value = '[Feng%Xiao%xref no email%1,          Rongqing%Sun%xref no email%1,          Wenbo%Sun%xref no email%1,          Dan%Xu%xref no email%1,          Lan%Lan%xref no email%1,          Huan%Li%xref no email%1,          Huan%Liu%xref no email%0,          Haibo%Xu%xref no email%1]'
$ws.Range("E32").Value = '[Dong%Ji%NULL%1,                 Dawei%Zhang%NULL%1,                 Jing%Xu%NULL%6,                 Zhu%Chen%NULL%1,                 Tieniu%Yang%NULL%1,                 Peng%Zhao%NULL%2,                 Guofeng%Chen%NULL%1,                 Gregory%Cheng%NULL%1,                 Yudong%Wang%NULL%1,                 Jingfeng%Bi%NULL%1,                 Lin%Tan%NULL%1,                 George%Lau%NULL%1,                 Enqiang%Qin%qeq2004@sina.com%1]'
$ws.Range("E33").Value = '[Fabrizio%Foieni%fabrizio.foieni@asst-valleolona.it%1,                 Girolamo%Sala%girolamo.sala@asst-valleolona.it%1,                 Jason Giuseppe%Mognarelli%giuseppe.mognarelli@unimi.it%1,                 Giulia%Suigo%NULL%1,                 Davide%Zampini%NULL%1,                 Matteo%Pistoia%NULL%1,                 Mariella%Ciola%NULL%1,                 Tommaso%Ciampani%NULL%1,                 Carolina%Ultori%NULL%1,                 Paolo%Ghiringhelli%NULL%1]'
$ws.Range("E34").Value = '[Alexander%Muacevic%NULL%5,                 John R%Adler%NULL%3,                 Sultan M%Kamran%NULL%2,                 Sultan M%Kamran%NULL%0,                 Zill-e-Humayun%Mirza%NULL%1,                 Hussain Abdul%Moeed%NULL%1,                 Arshad%Naseem%NULL%1,                 Maryam%Hussain%NULL%1,                 Imran%Fazal%NULL%1,                 Farrukh%Saeed%NULL%1,                 Wasim%Alamgir%NULL%1,                 Salman%Saleem%NULL%1,                 Sidra%Riaz%NULL%1]'
$ws.Range("E36").Value = '[Jie%Liu%xref no email%1,          Zilong%Liu%xref no email%1,          Weipeng%Jiang%xref no email%1,          Jian%Wang%xref no email%4,          Mengchan%Zhu%xref no email%1,          Juan%Song%xref no email%0,          Xiaoyue%Wang%xref no email%1,          Ying%Su%xref no email%1,          Guiling%Xiang%xref no email%1,          Maosong%Ye%xref no email%1,          Jiamin%Li%xref no email%1,          Yong%Zhang%xref no email%1,          Qinjun%Shen%xref no email%1,          Zhuozhe%Li%xref no email%1,          Danwei%Yao%xref no email%1,          Yuanlin%Song%xref no email%0,          Kaihuan%Yu%xref no email%1,          Zhe%Luo%xref no email%1,          Ling%Ye%xref no email%1]'
$ws.Range("E37").Value = '[Leonardo%Lorente%NULL%1,                 María M.%Martín%NULL%1,                 Mónica%Argueso%NULL%1,                 Jordi%Solé-Violán%NULL%1,                 Alina%Perez%NULL%1,                 José Alberto%Marcos Y Ramos%NULL%1,                 Luis%Ramos-Gómez%NULL%1,                 Sergio%López%NULL%1,                 Andrés%Franco%NULL%1,                 Agustín F.%González-Rivero%NULL%1,                 María%Martín%NULL%1,                 Verónica%Gonzalez%NULL%1,                 Julia%Alcoba-Flórez%NULL%1,                 Miguel Ángel%Rodriguez%NULL%1,                 Marta%Riaño-Ruiz%NULL%1,                 Juan%Guillermo O Campo%NULL%1,                 Lourdes%González%NULL%1,                 Tamara%Cantera%NULL%1,                 Raquel%Ortiz-López%NULL%1,                 Nazario%Ojeda%NULL%1,                 Aurelio%Rodríguez-Pérez%NULL%1,                 Casimira%Domínguez%NULL%1,                 Alejandro%Jiménez%NULL%1]'
$ws.Range("E38").Value = '[Gunther%Eysenbach%NULL%7,                 Rita%Kukafka%NULL%1,                 Arriel%Benis%NULL%2,                 Arriel%Benis%NULL%0,                 Jinfeng%Li%NULL%1,                 Pan%Pan%NULL%2,                 Pan%Pan%NULL%0,                 Yichao%Li%NULL%2,                 Yichao%Li%NULL%0,                 Yongjiu%Xiao%NULL%2,                 Yongjiu%Xiao%NULL%0,                 Bingchao%Han%NULL%2,                 Bingchao%Han%NULL%0,                 Longxiang%Su%NULL%2,                 Longxiang%Su%NULL%0,                 Mingliang%Su%NULL%2,                 Mingliang%Su%NULL%0,                 Yansheng%Li%NULL%2,                 Yansheng%Li%NULL%0,                 Siqi%Zhang%NULL%2,                 Siqi%Zhang%NULL%0,                 Dapeng%Jiang%NULL%2,                 Dapeng%Jiang%NULL%0,                 Xia%Chen%NULL%2,                 Xia%Chen%NULL%0,                 Fuquan%Zhou%NULL%2,                 Fuquan%Zhou%NULL%0,                 Ling%Ma%NULL%2,                 Ling%Ma%NULL%0,                 Pengtao%Bao%NULL%2,                 Pengtao%Bao%NULL%0,                 Lixin%Xie%xielx301@126.com%2,                 Lixin%Xie%xielx301@126.com%0]'
$ws.Range("E39").Value = '[Sonsoles%Salto-Alejandre%NULL%2,                 Cristina%Roca-Oporto%NULL%2,                 Guillermo%Martín-Gutiérrez%NULL%2,                 María Dolores%Avilés%NULL%2,                 Carmen%Gómez-González%NULL%2,                 María Dolores%Navarro-Amuedo%NULL%2,                 Julia%Praena-Segovia%NULL%2,                 José%Molina%NULL%2,                 María%Paniagua-García%NULL%2,                 Horacio%García-Delgado%NULL%2,                 Antonio%Domínguez-Petit%NULL%2,                 Jerónimo%Pachón%NULL%2,                 José Miguel%Cisneros%NULL%2]'
$ws.Range("E40").Value = '[Justin J.%Turcotte%NULL%1,                 Barry R.%Meisenberg%NULL%2,                 Barry R.%Meisenberg%NULL%0,                 James H.%MacDonald%NULL%1,                 Nandakumar%Menon%NULL%1,                 Marcia B.%Fowler%NULL%1,                 Michaline%West%NULL%1,                 Jane%Rhule%NULL%1,                 Sadaf S.%Qureshi%NULL%1,                 Eileen B.%MacDonald%NULL%1,                 Yu Ru%Kou%NULL%2,                 Yu Ru%Kou%NULL%0]'
$ws.Range("E41").Value = '[Gunther%Eysenbach%NULL%0,                 Sarah%Poole%NULL%2,                 Sarah%Poole%NULL%0,                 Meeta%Pradhan%NULL%1,                 Akhil%Vaid%NULL%3,                 Akhil%Vaid%NULL%0,                 Sulaiman%Somani%NULL%3,                 Sulaiman%Somani%NULL%0,                 Adam J%Russak%NULL%2,                 Adam J%Russak%NULL%0,                 Jessica K%De Freitas%NULL%2,                 Jessica K%De Freitas%NULL%0,                 Fayzan F%Chaudhry%NULL%2,                 Fayzan F%Chaudhry%NULL%0,                 Ishan%Paranjpe%NULL%3,                 Ishan%Paranjpe%NULL%0,                 Kipp W%Johnson%NULL%2,                 Kipp W%Johnson%NULL%0,                 Samuel J%Lee%NULL%2,                 Samuel J%Lee%NULL%0,                 Riccardo%Miotto%NULL%2,                 Riccardo%Miotto%NULL%0,                 Felix%Richter%NULL%3,                 Felix%Richter%NULL%0,                 Shan%Zhao%NULL%3,                 Shan%Zhao%NULL%0,                 Noam D%Beckmann%NULL%2,                 Noam D%Beckmann%NULL%0,                 Nidhi%Naik%NULL%2,                 Nidhi%Naik%NULL%0,                 Arash%Kia%NULL%2,                 Arash%Kia%NULL%0,                 Prem%Timsina%NULL%2,                 Prem%Timsina%NULL%0,                 Anuradha%Lala%NULL%2,                 Anuradha%Lala%NULL%0,                 Manish%Paranjpe%NULL%2,                 Manish%Paranjpe%NULL%0,                 Eddye%Golden%NULL%2,                 Eddye%Golden%NULL%0,                 Matteo%Danieletto%NULL%2,                 Matteo%Danieletto%NULL%0,                 Manbir%Singh%NULL%2,                 Manbir%Singh%NULL%0,                 Dara%Meyer%NULL%2,                 Dara%Meyer%NULL%0,                 Paul F%O''Reilly%NULL%2,                 Paul F%O''Reilly%NULL%0,                 Laura%Huckins%NULL%2,                 Laura%Huckins%NULL%0,                 Patricia%Kovatch%NULL%2,                 Patricia%Kovatch%NULL%0,                 Joseph%Finkelstein%NULL%2,                 Joseph%Finkelstein%NULL%0,                 Robert M.%Freeman%NULL%2,                 Robert M.%Freeman%NULL%0,                 Edgar%Argulian%NULL%2,                 Edgar%Argulian%NULL%0,                 Andrew%Kasarskis%NULL%2,                 Andrew%Kasarskis%NULL%0,                 Bethany%Percha%NULL%2,                 Bethany%Percha%NULL%0,                 Judith A%Aberg%NULL%2,                 Judith A%Aberg%NULL%0,                 Emilia%Bagiella%NULL%3,                 Emilia%Bagiella%NULL%0,                 Carol R%Horowitz%NULL%2,                 Carol R%Horowitz%NULL%0,                 Barbara%Murphy%NULL%2,                 Barbara%Murphy%NULL%0,                 Eric J%Nestler%NULL%2,                 Eric J%Nestler%NULL%0,                 Eric E%Schadt%NULL%2,                 Eric E%Schadt%NULL%0,                 Judy H%Cho%NULL%2,                 Judy H%Cho%NULL%0,                 Carlos%Cordon-Cardo%NULL%2,                 Carlos%Cordon-Cardo%NULL%0,                 Valentin%Fuster%NULL%3,                 Valentin%Fuster%NULL%0,                 Dennis S%Charney%NULL%2,                 Dennis S%Charney%NULL%0,                 David L%Reich%NULL%2,                 David L%Reich%NULL%0,                 Erwin P%Bottinger%NULL%2,                 Erwin P%Bottinger%NULL%0,                 Matthew A%Levin%NULL%2,                 Matthew A%Levin%NULL%0,                 Jagat%Narula%NULL%3,                 Jagat%Narula%NULL%0,                 Zahi A%Fayad%NULL%2,                 Zahi A%Fayad%NULL%0,                 Allan C%Just%NULL%2,                 Allan C%Just%NULL%0,                 Alexander W%Charney%NULL%2,                 Alexander W%Charney%NULL%0,                 Girish N%Nadkarni%NULL%2,                 Girish N%Nadkarni%NULL%0,                 Benjamin S%Glicksberg%benjamin.glicksberg@mssm.edu%2,                 Benjamin S%Glicksberg%benjamin.glicksberg@mssm.edu%0]'
$ws.Range("E42").Value = '[Zirun%Zhao%NULL%2,                 Anne%Chen%NULL%2,                 Wei%Hou%NULL%3,                 James M.%Graham%NULL%1,                 Haifang%Li%NULL%2,                 Paul S.%Richman%NULL%1,                 Henry C.%Thode%NULL%1,                 Adam J.%Singer%NULL%1,                 Tim Q.%Duong%NULL%1,                 Muhammad%Adrish%NULL%9,                 Muhammad%Adrish%NULL%0,                 Muhammad%Adrish%NULL%0,                 Muhammad%Adrish%NULL%0,                 Muhammad%Adrish%NULL%0,                 Muhammad%Adrish%NULL%0,                 Muhammad%Adrish%NULL%0,                 Muhammad%Adrish%NULL%0,                 Muhammad%Adrish%NULL%0]'
$ws.Range("E43").Value = '[Bin%Zhang%NULL%1,                 Qin%Liu%NULL%2,                 Xiao%Zhang%NULL%2,                 Shuyi%Liu%NULL%1,                 Weiqi%Chen%NULL%1,                 Jingjing%You%NULL%1,                 Qiuying%Chen%NULL%1,                 Minmin%Li%NULL%1,                 Zhuozhi%Chen%NULL%1,                 Luyan%Chen%NULL%1,                 Lv%Chen%NULL%1,                 Yuhao%Dong%NULL%1,                 Qingsi%Zeng%NULL%1,                 Shuixing%Zhang%NULL%1]'
$ws.Range("E44").Value = '[Jing%Yu%NULL%1,                 Lei%Nie%NULL%1,                 Dongde%Wu%NULL%1,                 Jian%Chen%NULL%2,                 Zhifeng%Yang%NULL%1,                 Ling%Zhang%NULL%3,                 Dongqing%Li%NULL%1,                 Xia%Zhou%NULL%1]'
$ws.Range("E45").Value = '[Binchen%Wang%xref no email%1,          Feiyang%Zhong%xref no email%1,          Hanfei%Zhang%xref no email%1,          Wenting%An%xref no email%1,          Meiyan%Liao%xref no email%1,          Yiyuan%Cao%xref no email%1]'
$ws.Range("E46").Value = '[Yelda%Varol%yeldavatansever@hotmail.com%1,                 Burcin%Hakoglu%NULL%2,                 Burcin%Hakoglu%NULL%0,                 Ali%Kadri Cirak%NULL%1,                 Gulru%Polat%NULL%1,                 Berna%Komurcuoglu%NULL%1,                 Berrin%Akkol%NULL%1,                 Cagri%Atasoy%NULL%1,                 Eda%Bayramic%NULL%1,                 Gunseli%Balci%NULL%1,                 Sena%Ataman%NULL%1,                 Sinem%Ermin%NULL%1,                 Enver%Yalniz%NULL%1,                 NULL%NULL%NULL%0]'
$ws.Range("E47").Value = '[Jing%Zhou%xref no email%1,          Lili%Huang%xref no email%1,          Jin%Chen%xref no email%1,          Xiaowei%Yuan%xref no email%1,          Qinhua%Shen%xref no email%1,          Su%Dong%xref no email%1,          Bei%Cheng%xref no email%1,          Tang-Meng%Guo%xref no email%1]'
$ws.Range("E48").Value = '[Zhihong%Weng%NULL%1,                 Qiaosen%Chen%NULL%2,                 Qiaosen%Chen%NULL%0,                 Sumeng%Li%NULL%1,                 Huadong%Li%NULL%1,                 Qian%Zhang%NULL%1,                 Sihong%Lu%NULL%1,                 Li%Wu%NULL%1,                 Leiqun%Xiong%NULL%1,                 Bobin%Mi%NULL%1,                 Di%Liu%NULL%3,                 Mengji%Lu%NULL%1,                 Dongliang%Yang%NULL%1,                 Hongbo%Jiang%hongbojiang3@163.com%1,                 Shaoping%Zheng%zhengspxx@126.com%1,                 Xin%Zheng%xin11@hotmail.com%1]'
$ws.Range("E49").Value = '[Junhong%Wang%NULL%1,                 Hua%Zhang%NULL%2,                 Rui%Qiao%NULL%1,                 Qinggang%Ge%NULL%1,                 Shuisheng%Zhang%NULL%1,                 Zongxuan%Zhao%NULL%1,                 Ci%Tian%NULL%1,                 Qingbian%Ma%NULL%2,                 Qingbian%Ma%NULL%0,                 Ning%Shen%NULL%1]'
$ws.Range("E50").Value = '[Stephen R%Knight%NULL%1,                 Antonia%Ho%NULL%2,                 Riinu%Pius%NULL%2,                 Iain%Buchan%NULL%1,                 Gail%Carson%NULL%2,                 Thomas M%Drake%NULL%1,                 Jake%Dunning%NULL%2,                 Cameron J%Fairfield%NULL%2,                 Carrol%Gamble%NULL%2,                 Christopher A%Green%NULL%2,                 Rishi%Gupta%NULL%1,                 Sophie%Halpin%NULL%2,                 Hayley E%Hardwick%NULL%1,                 Karl A%Holden%NULL%1,                 Peter W%Horby%NULL%2,                 Clare%Jackson%NULL%2,                 Kenneth A%Mclean%NULL%2,                 Laura%Merson%NULL%2,                 Jonathan S%Nguyen-Van-Tam%NULL%1,                 Lisa%Norman%NULL%2,                 Mahdad%Noursadeghi%NULL%2,                 Piero L%Olliaro%NULL%1,                 Mark G%Pritchard%NULL%1,                 Clark D%Russell%NULL%2,                 Catherine A%Shaw%NULL%2,                 Aziz%Sheikh%NULL%2,                 Tom%Solomon%NULL%2,                 Cathie%Sudlow%NULL%1,                 Olivia V%Swann%NULL%1,                 Lance CW%Turtle%NULL%2,                 Peter JM%Openshaw%NULL%2,                 J Kenneth%Baillie%NULL%2,                 Malcolm G%Semple%NULL%3,                 Annemarie B%Docherty%NULL%3,                 Annemarie B%Docherty%NULL%0,                 Ewen M%Harrison%NULL%3,                 NULL%NULL%NULL%0,                 NULL%NULL%NULL%0,                 J Kenneth%Baillie%NULL%0,                 Malcolm G%Semple%NULL%0,                 Peter JM%Openshaw%NULL%0,                 Gail%Carson%NULL%0,                 Beatrice%Alex%NULL%1,                 Benjamin%Bach%NULL%1,                 Wendy S%Barclay%NULL%1,                 Debby%Bogaert%NULL%1,                 Meera%Chand%NULL%1,                 Graham S%Cooke%NULL%1,                 Annemarie B%Docherty%NULL%0,                 Jake%Dunning%NULL%0,                 Ana%da Silva Filipe%NULL%1,                 Tom%Fletcher%NULL%1,                 Christopher A%Green%NULL%0,                 Ewen M%Harrison%NULL%0,                 Julian A%Hiscox%NULL%1,                 Antonia Ying Wai%Ho%NULL%1,                 Peter W%Horby%NULL%0,                 Samreen%Ijaz%NULL%1,                 Saye%Khoo%NULL%1,                 Paul%Klenerman%NULL%1,                 Andrew%Law%NULL%2,                 Wei Shen%Lim%NULL%1,                 Alexander J%Mentzer%NULL%1,                 Laura%Merson%NULL%0,                 Alison M%Meynert%NULL%1,                 Mahdad%Noursadeghi%NULL%0,                 Shona C%Moore%NULL%2,                 Massimo%Palmarini%NULL%1,                 William A%Paxton%NULL%1,                 Georgios%Pollakis%NULL%1,                 Nicholas%Price%NULL%1,                 Andrew%Rambaut%NULL%1,                 David L%Robertson%NULL%1,                 Clark D%Russell%NULL%0,                 Vanessa%Sancho-Shimizu%NULL%1,                 Janet T%Scott%NULL%1,                 Louise%Sigfrid%NULL%1,                 Tom%Solomon%NULL%0,                 Shiranee%Sriskandan%NULL%1,                 David%Stuart%NULL%1,                 Charlotte%Summers%NULL%1,                 Richard S%Tedder%NULL%1,                 Emma C%Thomson%NULL%1,                 Ryan S%Thwaites%NULL%1,                 Lance CW%Turtle%NULL%0,                 Maria%Zambon%NULL%1,                 Hayley%Hardwick%NULL%1,                 Chloe%Donohue%NULL%1,                 Jane%Ewins%NULL%1,                 Wilna%Oosthuyzen%NULL%1,                 Fiona%Griffiths%NULL%1,                 Lisa%Norman%NULL%0,                 Riinu%Pius%NULL%0,                 Tom M%Drake%NULL%1,                 Cameron J%Fairfield%NULL%0,                 Stephen%Knight%NULL%1,                 Kenneth A%Mclean%NULL%0,                 Derek%Murphy%NULL%1,                 Catherine A%Shaw%NULL%0,                 Jo%Dalton%NULL%1,                 Michelle%Girvan%NULL%1,                 Egle%Saviciute%NULL%1,                 Stephanie%Roberts%NULL%1,                 Janet%Harrison%NULL%1,                 Laura%Marsh%NULL%1,                 Marie%Connor%NULL%1,                 Sophie%Halpin%NULL%0,                 Clare%Jackson%NULL%0,                 Carrol%Gamble%NULL%0,                 Gary%Leeming%NULL%1,                 Andrew%Law%NULL%0,                 Ross%Hendry%NULL%1,                 James%Scott-Brown%NULL%1,                 William%Greenhalf%NULL%1,                 Victoria%Shaw%NULL%1,                 Sarah%McDonald%NULL%2,                 Katie A%Ahmed%NULL%1,                 Jane A%Armstrong%NULL%1,                 Milton%Ashworth%NULL%1,                 Innocent G%Asiimwe%NULL%1,                 Siddharth%Bakshi%NULL%1,                 Samantha L%Barlow%NULL%1,                 Laura%Booth%NULL%1,                 Benjamin%Brennan%NULL%1,                 Katie%Bullock%NULL%1,                 Benjamin WA%Catterall%NULL%1,                 Jordan J%Clark%NULL%1,                 Emily A%Clarke%NULL%1,                 Sarah%Cole%NULL%1,                 Louise%Cooper%NULL%1,                 Helen%Cox%NULL%1,                 Christopher%Davis%NULL%1,                 Oslem%Dincarslan%NULL%1,                 Chris%Dunn%NULL%1,                 Philip%Dyer%NULL%1,                 Angela%Elliott%NULL%1,                 Anthony%Evans%NULL%1,                 Lewis WS%Fisher%NULL%1,                 Terry%Foster%NULL%1,                 Isabel%Garcia-Dorival%NULL%1,                 Willliam%Greenhalf%NULL%1,                 Philip%Gunning%NULL%1,                 Catherine%Hartley%NULL%1,                 Antonia%Ho%NULL%0,                 Rebecca L%Jensen%NULL%1,                 Christopher B%Jones%NULL%1,                 Trevor R%Jones%NULL%1,                 Shadia%Khandaker%NULL%1,                 Katharine%King%NULL%1,                 Robyn T%Kiy%NULL%1,                 Chrysa%Koukorava%NULL%1,                 Annette%Lake%NULL%1,                 Suzannah%Lant%NULL%1,                 Diane%Latawiec%NULL%1,                 L%Lavelle-Langham%NULL%1,                 Daniella%Lefteri%NULL%1,                 Lauren%Lett%NULL%1,                 Lucia A%Livoti%NULL%1,                 Maria%Mancini%NULL%1,                 Sarah%McDonald%NULL%0,                 Laurence%McEvoy%NULL%1,                 John%McLauchlan%NULL%1,                 Soeren%Metelmann%NULL%1,                 Nahida S%Miah%NULL%1,                 Joanna%Middleton%NULL%1,                 Joyce%Mitchell%NULL%1,                 Shona C%Moore%NULL%0,                 Ellen G%Murphy%NULL%1,                 Rebekah%Penrice-Randal%NULL%1,                 Jack%Pilgrim%NULL%1,                 Tessa%Prince%NULL%1,                 Will%Reynolds%NULL%1,                 P Matthew%Ridley%NULL%1,                 Debby%Sales%NULL%1,                 Victoria E%Shaw%NULL%1,                 Rebecca K%Shears%NULL%1,                 Benjamin%Small%NULL%1,                 Krishanthi S%Subramaniam%NULL%1,                 Agnieska%Szemiel%NULL%1,                 Aislynn%Taggart%NULL%1,                 Jolanta%Tanianis-Hughes%NULL%1,                 Jordan%Thomas%NULL%1,                 Erwan%Trochu%NULL%1,                 Libby%van Tonder%NULL%1,                 Eve%Wilcock%NULL%1,                 J Eunice%Zhang%NULL%1,                 Kayode%Adeniji%NULL%1,                 Daniel%Agranoff%NULL%1,                 Ken%Agwuh%NULL%1,                 Dhiraj%Ail%NULL%1,                 Ana%Alegria%NULL%1,                 Brian%Angus%NULL%1,                 Abdul%Ashish%NULL%1,                 Dougal%Atkinson%NULL%1,                 Shahedal%Bari%NULL%1,                 Gavin%Barlow%NULL%1,                 Stella%Barnass%NULL%1,                 Nicholas%Barrett%NULL%2,                 Christopher%Bassford%NULL%1,                 David%Baxter%NULL%1,                 Michael%Beadsworth%NULL%1,                 Jolanta%Bernatoniene%NULL%1,                 John%Berridge%NULL%1,                 Nicola%Best%NULL%1,                 Pieter%Bothma%NULL%1,                 David%Brealey%NULL%1,                 Robin%Brittain-Long%NULL%1,                 Naomi%Bulteel%NULL%1,                 Tom%Burden%NULL%1,                 Andrew%Burtenshaw%NULL%1,                 Vikki%Caruth%NULL%1,                 David%Chadwick%NULL%1,                 Duncan%Chambler%NULL%1,                 Nigel%Chee%NULL%1,                 Jenny%Child%NULL%1,                 Srikanth%Chukkambotla%NULL%1,                 Tom%Clark%NULL%1,                 Paul%Collini%NULL%1,                 Catherine%Cosgrove%NULL%1,                 Jason%Cupitt%NULL%1,                 Maria-Teresa%Cutino-Moguel%NULL%1,                 Paul%Dark%NULL%1,                 Chris%Dawson%NULL%1,                 Samir%Dervisevic%NULL%1,                 Phil%Donnison%NULL%1,                 Sam%Douthwaite%NULL%1,                 Ingrid%DuRand%NULL%1,                 Ahilanadan%Dushianthan%NULL%1,                 Tristan%Dyer%NULL%1,                 Cariad%Evans%NULL%1,                 Chi%Eziefula%NULL%1,                 Chrisopher%Fegan%NULL%1,                 Adam%Finn%NULL%1,                 Duncan%Fullerton%NULL%1,                 Sanjeev%Garg%NULL%2,                 Sanjeev%Garg%NULL%0,                 Atul%Garg%NULL%1,                 Jo%Godden%NULL%1,                 Arthur%Goldsmith%NULL%1,                 Clive%Graham%NULL%1,                 Elaine%Hardy%NULL%1,                 Stuart%Hartshorn%NULL%1,                 Daniel%Harvey%NULL%1,                 Peter%Havalda%NULL%1,                 Daniel B%Hawcutt%NULL%1,                 Maria%Hobrok%NULL%1,                 Luke%Hodgson%NULL%1,                 Anita%Holme%NULL%1,                 Anil%Hormis%NULL%1,                 Michael%Jacobs%NULL%1,                 Susan%Jain%NULL%1,                 Paul%Jennings%NULL%1,                 Agilan%Kaliappan%NULL%1,                 Vidya%Kasipandian%NULL%1,                 Stephen%Kegg%NULL%1,                 Michael%Kelsey%NULL%1,                 Jason%Kendall%NULL%1,                 Caroline%Kerrison%NULL%1,                 Ian%Kerslake%NULL%1,                 Oliver%Koch%NULL%2,                 Gouri%Koduri%NULL%1,                 George%Koshy%NULL%1,                 Shondipon%Laha%NULL%1,                 Susan%Larkin%NULL%1,                 Tamas%Leiner%NULL%1,                 Patrick%Lillie%NULL%1,                 James%Limb%NULL%1,                 Vanessa%Linnett%NULL%1,                 Jeff%Little%NULL%1,                 Michael%MacMahon%NULL%1,                 Emily%MacNaughton%NULL%1,                 Ravish%Mankregod%NULL%1,                 Huw%Masson%NULL%1,                 Elijah%Matovu%NULL%1,                 Katherine%McCullough%NULL%1,                 Ruth%McEwen%NULL%1,                 Manjula%Meda%NULL%1,                 Gary%Mills%NULL%1,                 Jane%Minton%NULL%1,                 Mariyam%Mirfenderesky%NULL%1,                 Kavya%Mohandas%NULL%1,                 Quen%Mok%NULL%1,                 James%Moon%NULL%1,                 Elinoor%Moore%NULL%1,                 Patrick%Morgan%NULL%1,                 Craig%Morris%NULL%1,                 Katherine%Mortimore%NULL%1,                 Samuel%Moses%NULL%1,                 Mbiye%Mpenge%NULL%1,                 Rohinton%Mulla%NULL%1,                 Michael%Murphy%NULL%1,                 Megan%Nagel%NULL%1,                 Thapas%Nagarajan%NULL%1,                 Mark%Nelson%NULL%1,                 Igor%Otahal%NULL%1,                 Mark%Pais%NULL%1,                 Selva%Panchatsharam%NULL%1,                 Hassan%Paraiso%NULL%1,                 Brij%Patel%NULL%1,                 Justin%Pepperell%NULL%1,                 Mark%Peters%NULL%1,                 Mandeep%Phull%NULL%1,                 Stefania%Pintus%NULL%1,                 Jagtur Singh%Pooni%NULL%1,                 Frank%Post%NULL%1,                 David%Price%NULL%1,                 Rachel%Prout%NULL%1,                 Nikolas%Rae%NULL%1,                 Henrik%Reschreiter%NULL%1,                 Tim%Reynolds%NULL%1,                 Neil%Richardson%NULL%1,                 Mark%Roberts%NULL%1,                 Devender%Roberts%NULL%1,                 Alistair%Rose%NULL%1,                 Guy%Rousseau%NULL%1,                 Brendan%Ryan%NULL%1,                 Taranprit%Saluja%NULL%1,                 Aarti%Shah%NULL%1,                 Prad%Shanmuga%NULL%1,                 Anil%Sharma%NULL%1,                 Anna%Shawcross%NULL%1,                 Jeremy%Sizer%NULL%1,                 Richard%Smith%NULL%1,                 Catherine%Snelson%NULL%1,                 Nick%Spittle%NULL%1,                 Nikki%Staines%NULL%1,                 Tom%Stambach%NULL%1,                 Richard%Stewart%NULL%1,                 Pradeep%Subudhi%NULL%1,                 Tamas%Szakmany%NULL%1,                 Kate%Tatham%NULL%1,                 Jo%Thomas%NULL%1,                 Chris%Thompson%NULL%1,                 Robert%Thompson%NULL%1,                 Ascanio%Tridente%NULL%1,                 Darell%Tupper-Carey%NULL%1,                 Mary%Twagira%NULL%1,                 Andrew%Ustianowski%NULL%1,                 Nick%Vallotton%NULL%1,                 Lisa%Vincent-Smith%NULL%1,                 Shico%Visuvanathan%NULL%1,                 Alan%Vuylsteke%NULL%1,                 Sam%Waddy%NULL%1,                 Rachel%Wake%NULL%1,                 Andrew%Walden%NULL%1,                 Ingeborg%Welters%NULL%1,                 Tony%Whitehouse%NULL%1,                 Paul%Whittaker%NULL%1,                 Ashley%Whittington%NULL%1,                 Meme%Wijesinghe%NULL%1,                 Martin%Williams%NULL%1,                 Lawrence%Wilson%NULL%1,                 Sarah%Wilson%NULL%1,                 Stephen%Winchester%NULL%1,                 Martin%Wiselka%NULL%1,                 Adam%Wolverson%NULL%1,                 Daniel G%Wooton%NULL%1,                 Andrew%Workman%NULL%1,                 Bryan%Yates%NULL%1,                 Peter%Young%NULL%1]'
$ws.Range("E51").Value = '[Adrian%Soto‐Mota%adrian.soto@dpag.ox.ac.uk%1,                 Braulio A.%Marfil‐Garza%NULL%1,                 Erick%Martínez Rodríguez%NULL%1,                 José Omar%Barreto Rodríguez%NULL%1,                 Alicia Estela%López Romo%NULL%1,                 Paolo%Alberti Minutti%NULL%1,                 Juan Vicente%Alejandre Loya%NULL%1,                 Félix Emmanuel%Pérez Talavera%NULL%1,                 Freddy José%Ávila Cervera%NULL%1,                 Adriana%Velazquez Burciaga%NULL%1,                 Oscar%Morado Aramburo%NULL%1,                 Luis Alberto%Piña Olguín%NULL%1,                 Adrian%Soto‐Rodríguez%NULL%1,                 Andrés%Castañeda Prado%NULL%1,                 Patricio%Santillán Doherty%NULL%1,                 Juan%O Galindo%NULL%1,                 Luis Alberto%Guízar García%NULL%1,                 Daniel%Hernández Gordillo%NULL%1,                 Juan%Gutiérrez Mejía%NULL%1]'
$ws.Range("E52").Value = '[Ahmed%Abdulaal%NULL%1,                 Aatish%Patel%NULL%1,                 Esmita%Charani%NULL%1,                 Sarah%Denny%NULL%1,                 Saleh A.%Alqahtani%NULL%1,                 Gary W.%Davies%NULL%1,                 Nabeela%Mughal%NULL%1,                 Luke S. P.%Moore%l.moore@imperial.ac.uk%1]'
$ws.Range("E53").Value = '[Chansik%An%NULL%1,                 Hyunsun%Lim%NULL%1,                 Dong-Wook%Kim%NULL%1,                 Jung Hyun%Chang%NULL%1,                 Yoon Jung%Choi%chris316@yuhs.ac%1,                 Seong Woo%Kim%NULL%1]'
$ws.Range("E54").Value = '[Hao%Chen%xref no email%1,          Rudong%Chen%xref no email%1,          Hongkuan%Yang%xref no email%1,          Junhong%Wang%xref no email%1,          Yuyang%Hou%xref no email%1,          Wei%Hu%xref no email%1,          Jiasheng%Yu%xref no email%1,          Hua%Li%xref no email%1]'
$ws.Range("E55").Value = '[Ruchong%Chen%NULL%6,                 Wenhua%Liang%NULL%3,                 Mei%Jiang%NULL%2,                 Weijie%Guan%NULL%3,                 Chen%Zhan%NULL%2,                 Tao%Wang%NULL%0,                 Chunli%Tang%NULL%3,                 Ling%Sang%NULL%3,                 Jiaxing%Liu%NULL%2,                 Zhengyi%Ni%NULL%2,                 Yu%Hu%NULL%0,                 Lei%Liu%NULL%0,                 Hong%Shan%NULL%5,                 Chunliang%Lei%NULL%2,                 Yixiang%Peng%NULL%2,                 Li%Wei%NULL%5,                 Yong%Liu%NULL%5,                 Yahua%Hu%NULL%2,                 Peng%Peng%NULL%7,                 Jianming%Wang%NULL%2,                 Jiyang%Liu%NULL%2,                 Zhong%Chen%NULL%5,                 Gang%Li%NULL%5,                 Zhijian%Zheng%NULL%2,                 Shaoqin%Qiu%NULL%2,                 Jie%Luo%NULL%5,                 Changjiang%Ye%NULL%2,                 Shaoyong%Zhu%NULL%2,                 Xiaoqing%Liu%NULL%2,                 Linling%Cheng%NULL%2,                 Feng%Ye%NULL%2,                 Jinping%Zheng%NULL%2,                 Nuofu%Zhang%NULL%2,                 Yimin%Li%NULL%2,                 Jianxing%He%NULL%2,                 Shiyue%Li%lishiyue@188.com%3,                 Nanshan%Zhong%NULL%4,                 NULL%NULL%NULL%0]'
$ws.Range("E56").Value = '[Anying%Cheng%NULL%3,                 Liu%Hu%NULL%1,                 Yiru%Wang%NULL%0,                 Luyan%Huang%NULL%1,                 Lingxi%Zhao%NULL%1,                 Congcong%Zhang%NULL%1,                 Xiyue%Liu%NULL%1,                 Ranran%Xu%NULL%1,                 Feng%Liu%NULL%0,                 Jinping%Li%NULL%1,                 Dawei%Ye%NULL%0,                 Tao%Wang%NULL%0,                 Yongman%Lv%lvyongman@126.com%0,                 Qingquan%Liu%qqliutj@163.com%1]'
$ws.Range("E57").Value = '[Ash K%Clift%NULL%1,                 Carol A C%Coupland%NULL%1,                 Ruth H%Keogh%NULL%1,                 Karla%Diaz-Ordaz%NULL%1,                 Elizabeth%Williamson%NULL%1,                 Ewen M%Harrison%NULL%0,                 Andrew%Hayward%NULL%1,                 Harry%Hemingway%NULL%1,                 Peter%Horby%NULL%1,                 Nisha%Mehta%NULL%1,                 Jonathan%Benger%NULL%1,                 Kamlesh%Khunti%NULL%0,                 David%Spiegelhalter%NULL%1,                 Aziz%Sheikh%NULL%0,                 Jonathan%Valabhji%NULL%0,                 Ronan A%Lyons%NULL%1,                 John%Robson%NULL%1,                 Malcolm G%Semple%NULL%0,                 Frank%Kee%NULL%1,                 Peter%Johnson%NULL%1,                 Susan%Jebb%NULL%1,                 Tony%Williams%NULL%1,                 Julia%Hippisley-Cox%NULL%1]'
$ws.Range("E58").Value = '[Yue%Gao%NULL%1,                 Guang-Yao%Cai%NULL%1,                 Wei%Fang%NULL%2,                 Hua-Yi%Li%NULL%1,                 Si-Yuan%Wang%NULL%2,                 Si-Yuan%Wang%NULL%0,                 Lingxi%Chen%NULL%1,                 Yang%Yu%NULL%1,                 Dan%Liu%NULL%0,                 Sen%Xu%NULL%2,                 Peng-Fei%Cui%NULL%1,                 Shao-Qing%Zeng%NULL%2,                 Shao-Qing%Zeng%NULL%0,                 Xin-Xia%Feng%NULL%1,                 Rui-Di%Yu%NULL%1,                 Ya%Wang%NULL%2,                 Yuan%Yuan%NULL%1,                 Xiao-Fei%Jiao%NULL%1,                 Jian-Hua%Chi%NULL%1,                 Jia-Hao%Liu%NULL%1,                 Ru-Yuan%Li%NULL%1,                 Xu%Zheng%NULL%1,                 Chun-Yan%Song%NULL%1,                 Ning%Jin%NULL%1,                 Wen-Jian%Gong%NULL%1,                 Xing-Yu%Liu%NULL%1,                 Lei%Huang%NULL%2,                 Xun%Tian%NULL%1,                 Lin%Li%NULL%1,                 Hui%Xing%NULL%1,                 Ding%Ma%NULL%1,                 Chun-Rui%Li%NULL%1,                 Fei%Ye%yeyuanbei@hotmail.com%1,                 Qing-Lei%Gao%qingleigao@hotmail.com%2,                 Qing-Lei%Gao%qingleigao@hotmail.com%0]'
$ws.Range("E59").Value = '[Xiaoxu%Ma%NULL%1,                 Ang%Li%NULL%1,                 Mengfan%Jiao%NULL%1,                 Qingmiao%Shi%NULL%1,                 Xiaocai%An%NULL%1,                 Yonghai%Feng%NULL%1,                 Lihua%Xing%NULL%1,                 Hongxia%Liang%NULL%1,                 Jiajun%Chen%NULL%1,                 Huiling%Li%NULL%1,                 Juan%Li%NULL%0,                 Zhigang%Ren%NULL%1,                 Ranran%Sun%NULL%1,                 Guangying%Cui%NULL%1,                 Yongjian%Zhou%NULL%1,                 Ming%Cheng%NULL%1,                 Pengfei%Jiao%NULL%1,                 Yu%Wang%NULL%2,                 Jiyuan%Xing%NULL%1,                 Shen%Shen%NULL%1,                 Qingxian%Zhang%NULL%1,                 Aiguo%Xu%NULL%1,                 Zujiang%Yu%NULL%1]'
$ws.Range("E60").Value = '[Xuedi%Ma%NULL%1,                 Michael%Ng%NULL%1,                 Shuang%Xu%NULL%1,                 Zhouming%Xu%NULL%1,                 Hui%Qiu%NULL%1,                 Yuwei%Liu%NULL%1,                 Jiayou%Lyu%NULL%1,                 Jiwen%You%NULL%1,                 Peng%Zhao%NULL%0,                 Shihao%Wang%NULL%1,                 Yunfei%Tang%NULL%1,                 Hao%Cui%NULL%1,                 Changxiao%Yu%NULL%1,                 Feng%Wang%NULL%5,                 Fei%Shao%NULL%1,                 Peng%Sun%NULL%1,                 Ziren%Tang%NULL%1]'
$ws.Range("E61").Value = '[Xiaojun%Ma%NULL%1,                 Huifang%Wang%NULL%1,                 Junwei%Huang%NULL%1,                 Yan%Geng%NULL%1,                 Shuqi%Jiang%NULL%1,                 Qiuping%Zhou%NULL%1,                 Xuan%Chen%NULL%1,                 Hongping%Hu%NULL%1,                 Weifeng%Li%NULL%1,                 Chengbin%Zhou%NULL%1,                 Xinglin%Gao%NULL%1,                 Na%Peng%pnatz@163.com%1,                 Yiyu%Deng%yiyudeng666@163.com%1]'
$ws.Range("E62").Value = '[Deng%Pan%NULL%1,                 Dandan%Cheng%NULL%1,                 Yiwei%Cao%NULL%1,                 Chuan%Hu%NULL%1,                 Fenglin%Zou%NULL%1,                 Wencheng%Yu%NULL%1,                 Tao%Xu%NULL%2]'
$ws.Range("E63").Value = '[Jung Gil%Park%NULL%1,                 Min Kyu%Kang%NULL%2,                 Min Kyu%Kang%NULL%0,                 Yu Rim%Lee%NULL%2,                 Yu Rim%Lee%NULL%0,                 Jeong Eun%Song%NULL%2,                 Jeong Eun%Song%NULL%0,                 Na Young%Kim%NULL%1,                 Young Oh%Kweon%NULL%1,                 Won Young%Tak%NULL%1,                 Se Young%Jang%NULL%1,                 Changhyeong%Lee%NULL%2,                 Changhyeong%Lee%NULL%0,                 Byung Seok%Kim%NULL%1,                 Jae Seok%Hwang%NULL%1,                 Byoung Kuk%Jang%NULL%1,                 Jinmok%Bae%NULL%1,                 Ji Yeon%Lee%NULL%0,                 Jeong Ill%Suh%NULL%1,                 Soo Young%Park%NULL%1,                 Woo Jin%Chung%NULL%2,                 Woo Jin%Chung%NULL%0,                 NULL%NULL%NULL%0]'
$ws.Range("E64").Value = '[Sonsoles%Salto-Alejandre%NULL%0,                 Cristina%Roca-Oporto%NULL%0,                 Guillermo%Martín-Gutiérrez%NULL%0,                 María Dolores%Avilés%NULL%0,                 Carmen%Gómez-González%NULL%0,                 María Dolores%Navarro-Amuedo%NULL%0,                 Julia%Praena-Segovia%NULL%0,                 José%Molina%NULL%0,                 María%Paniagua-García%NULL%0,                 Horacio%García-Delgado%NULL%0,                 Antonio%Domínguez-Petit%NULL%0,                 Jerónimo%Pachón%NULL%0,                 José Miguel%Cisneros%NULL%0]'
$ws.Range("E65").Value = '[Arenas%Joaqu\u00edn%coreGivesNoEmail%1,              Calvo%Boyero Fernando%coreGivesNoEmail%1,              Castillo%Garc\u00eda Adri\u00e1n%coreGivesNoEmail%1,              COVID-19%Hospital \u201912 Octubre\u2019 Clinical Biochemisty Study Group%coreGivesNoEmail%1,              Cueto%Felgueroso Cecilia%coreGivesNoEmail%1,              Luc\u00eda%Mulas Alejandro%coreGivesNoEmail%1,              L\u00f3pez%Jim\u00e9nez Ana%coreGivesNoEmail%1,              Mart\u00edn%Casanueva Miguel \u00c1ngel%coreGivesNoEmail%1,              Santos%Lozano Alejandro%coreGivesNoEmail%1,              Valenzuela%Ruiz Pedro Luis%coreGivesNoEmail%1]'
